# Generate Report for Handoff
# Update the localization status report: the zh-cn item has moved from
# "In Translation" to "Ready for handoff", with refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status column updates (zh-cn and de-de status, both shown on Overview
# as well as on each language-specific sheet)
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# Latest handoff date/time updates
$overview.Range("D2").Value = "2016-03-24 06:40:59"
$zhcn.Range("E2").Value     = "2016-03-24 06:40:55"
$dede.Range("E2").Value     = "2016-03-24 06:40:59"
